# feat: add single and multi corrector
# Applies corrected reference labels / ranges to the CBC report sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number need to stay text
# (matches the source data, which stores every cell as a string) so we
# pre-format them as Text before writing the value.
$textCells = @("B5", "B6", "B7", "B13", "B14", "B18", "C29")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("A2").Value = "中性细胞数"
$ws.Range("C2").Value = "4-10"

$ws.Range("A3").Value = "中性细胞数"
$ws.Range("C3").Value = "3.5-5"

$ws.Range("C4").Value = "110-150"

$ws.Range("B5").Value = "0.38"
$ws.Range("C5").Value = "0.37-0.47"

$ws.Range("B6").Value = "90.3"
$ws.Range("C6").Value = "80-100"

$ws.Range("B7").Value = "29.9"
$ws.Range("C7").Value = "27-34"

$ws.Range("A8").Value = "RDW-CV"
$ws.Range("C8").Value = "320-360"

$ws.Range("A9").Value = "RDW-CV"

$ws.Range("C10").Value = "11.6-14.6"

$ws.Range("C11").Value = "100-300"

$ws.Range("C12").Value = "0.11-0.28"

$ws.Range("A13").Value = "血小板分布宽度"
$ws.Range("B13").Value = "9.9"
$ws.Range("C13").Value = "9-17"

$ws.Range("A14").Value = "平均血小板体积"
$ws.Range("B14").Value = "9.7"
$ws.Range("C14").Value = "6.5-11"

$ws.Range("A15").Value = "中性粒细胞数"
$ws.Range("C15").Value = "×"

$ws.Range("A16").Value = "淋巴细胞数"
$ws.Range("C16").Value = "X"

$ws.Range("A17").Value = "单核细胞数"
$ws.Range("C17").Value = "X"

$ws.Range("A18").Value = "嗜酸性粒细胞数"
$ws.Range("B18").Value = "0.02"
$ws.Range("C18").Value = "×"

$ws.Range("A19").Value = "RDW-CV"

$ws.Range("A20").Value = "RDW-CV"
$ws.Range("C20").Value = "2013-03"

$ws.Range("A21").Value = "中性粒细胞数"
$ws.Range("C21").Value = "×"

$ws.Range("A22").Value = "中性粒细胞绝对值"
$ws.Range("C22").Value = "×10^91"

$ws.Range("A23").Value = "中性粒细胞绝对值"

$ws.Range("A24").Value = "淋巴细胞绝对值"

$ws.Range("A25").Value = "单核细胞绝对值"

$ws.Range("A26").Value = "嗜酸性粒细胞绝对值"

$ws.Range("A27").Value = "嗜碱性粒细胞绝对值"

$ws.Range("C29").Value = "1003"

$ws.Range("A34").Value = "RDW-CV"
